# Sprint 3 kickoff: add a new "CertificateTemplatesPath" config row to the
# LiveConfig Sheet1 table (directly above "CountryLookUpPath"), grow the
# Table1 range to match, and restore the view state (active sheet/selection)
# used while making the edit.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Insert the new parameter row right before "CountryLookUpPath" (row 26) ---
$ws1.Rows("26:26").Insert()
$ws1.Rows("26:26").RowHeight = 48.75

$ws1.Range("A26").Value = "CertificateTemplatesPath"
$ws1.Range("B26").Value = "\\EARTH.GSI.GOV.UK\USER\SHARED\Agency\CoFS for G drive\RobotDocuments\Robot Certificate Templates\"
$ws1.Range("C26").Value = "Folder path for the robot's version of certificate templates"

# --- Grow Table1 (and its filter range) to include the newly inserted row ---
$lo = $ws1.ListObjects.Item("Table1")
$lo.Resize($ws1.Range("A1:C35"))

# --- Restore sheet/selection state: Sheet1 active, scrolled near the new
#     row, with B27 (the shifted "Value" cell for CountryLookUpPath) selected ---
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("B27").Select()
